# Auto-generated edit script applying 187 cell-value updates
# (scheduled-runner refresh of market price / profit columns H:N)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3542.1538
$ws.Range("I17").Value = 747.2222
$ws.Range("J17").Value = 4380.6333
$ws.Range("K17").Value = 2241.6666
$ws.Range("L17").Value = 13141.8999
$ws.Range("M17").Value = -2073.6666
$ws.Range("N17").Value = -13477.8999
$ws.Range("H70").Value = 75528.5
$ws.Range("J70").Value = 4900.3
$ws.Range("L70").Value = 14700.9
$ws.Range("N70").Value = -15240.9
$ws.Range("H73").Value = 75528.5
$ws.Range("J73").Value = 4900.3
$ws.Range("L73").Value = 14700.9
$ws.Range("N73").Value = -16572.9
$ws.Range("H132").Value = 1180.5084
$ws.Range("I132").Value = 1181.193
$ws.Range("K132").Value = 3543.579
$ws.Range("M132").Value = -1013.579
$ws.Range("H135").Value = 1831.8
$ws.Range("I135").Value = 1259.8572
$ws.Range("K135").Value = 11338.7148
$ws.Range("M135").Value = -8803.7148
$ws.Range("H137").Value = 28574574
$ws.Range("I137").Value = 58825970
$ws.Range("K137").Value = 176477910
$ws.Range("M137").Value = -176475360
$ws.Range("H138").Value = 3486.2058
$ws.Range("I138").Value = 2027
$ws.Range("K138").Value = 6081
$ws.Range("M138").Value = -941
$ws.Range("H141").Value = 740.25
$ws.Range("I141").Value = 740.25
$ws.Range("K141").Value = 2220.75
$ws.Range("M141").Value = 2959.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5012.6665
$ws.Range("I61").Value = 5103.829
$ws.Range("K61").Value = 5103.829
$ws.Range("M61").Value = -4891.829
$ws.Range("H110").Value = 1740.6111
$ws.Range("I110").Value = 726.125
$ws.Range("K110").Value = 726.125
$ws.Range("M110").Value = 1318.875
$ws.Range("H132").Value = 7070.6787
$ws.Range("I132").Value = 2854.3333
$ws.Range("J132").Value = 14660.1
$ws.Range("K132").Value = 8562.999899999999
$ws.Range("L132").Value = 43980.3
$ws.Range("M132").Value = -6032.999899999999
$ws.Range("N132").Value = -49040.3
$ws.Range("H136").Value = 5012.6665
$ws.Range("I136").Value = 5103.829
$ws.Range("K136").Value = 15311.487
$ws.Range("M136").Value = -12761.487

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 26991.857
$ws.Range("I75").Value = 11899.8
$ws.Range("K75").Value = 11899.8
$ws.Range("M75").Value = -10963.8
$ws.Range("H78").Value = 26991.857
$ws.Range("I78").Value = 11899.8
$ws.Range("K78").Value = 35699.39999999999
$ws.Range("M78").Value = -31019.39999999999
$ws.Range("H134").Value = 1810.4117
$ws.Range("I134").Value = 1290.3334
$ws.Range("K134").Value = 3871.0002
$ws.Range("M134").Value = -1336.0002
$ws.Range("H140").Value = 58702.832
$ws.Range("J140").Value = 58702.832
$ws.Range("L140").Value = 58702.832
$ws.Range("N140").Value = -69062.83199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2400.64
$ws.Range("I16").Value = 1492.4706
$ws.Range("K16").Value = 1492.4706
$ws.Range("M16").Value = -1205.4706
$ws.Range("H31").Value = 64519.59
$ws.Range("J31").Value = 78130.92999999999
$ws.Range("L31").Value = 78130.92999999999
$ws.Range("N31").Value = -78720.92999999999
$ws.Range("H34").Value = 64519.59
$ws.Range("J34").Value = 78130.92999999999
$ws.Range("L34").Value = 78130.92999999999
$ws.Range("N34").Value = -78534.92999999999
$ws.Range("H58").Value = 2980.9285
$ws.Range("I58").Value = 1176.75
$ws.Range("K58").Value = 1176.75
$ws.Range("M58").Value = -973.75
$ws.Range("H113").Value = 2400.64
$ws.Range("I113").Value = 1492.4706
$ws.Range("K113").Value = 1492.4706
$ws.Range("M113").Value = 677.5293999999999
$ws.Range("H136").Value = 2980.9285
$ws.Range("I136").Value = 1176.75
$ws.Range("K136").Value = 3530.25
$ws.Range("M136").Value = -980.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 27.5
$ws.Range("I13").Value = 27.5
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 82.5
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 85.5
$ws.Range("N13").ClearContents()
$ws.Range("H33").Value = 666719.8
$ws.Range("I33").Value = 1000040.7
$ws.Range("J33").Value = 78
$ws.Range("K33").Value = 6000244.199999999
$ws.Range("L33").Value = 468
$ws.Range("M33").Value = -5999961.199999999
$ws.Range("N33").Value = -1034
$ws.Range("H39").Value = 1900
$ws.Range("I39").Value = 1537.5
$ws.Range("J39").Value = 4800
$ws.Range("K39").Value = 4612.5
$ws.Range("L39").Value = 14400
$ws.Range("M39").Value = -4318.5
$ws.Range("N39").Value = -14988
$ws.Range("H113").Value = 90910470
$ws.Range("I113").Value = 1732.6666
$ws.Range("J113").Value = 125001250
$ws.Range("K113").Value = 5197.9998
$ws.Range("L113").Value = 375003750
$ws.Range("M113").Value = -3027.9998
$ws.Range("N113").Value = -375008090
$ws.Range("H132").Value = 3455.75
$ws.Range("J132").Value = 3963.55
$ws.Range("L132").Value = 35671.95
$ws.Range("N132").Value = -40731.95

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 308151.62
$ws.Range("I132").Value = 339038.25
$ws.Range("K132").Value = 1017114.75
$ws.Range("M132").Value = -1014584.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7594.5
$ws.Range("I7").Value = 4434.7896
$ws.Range("J7").Value = 16170.857
$ws.Range("K7").Value = 4434.7896
$ws.Range("L7").Value = 16170.857
$ws.Range("M7").Value = -4322.7896
$ws.Range("N7").Value = -16394.857
$ws.Range("H40").Value = 21890
$ws.Range("J40").Value = 19005
$ws.Range("L40").Value = 19005
$ws.Range("N40").Value = -19277
$ws.Range("H126").Value = 7594.5
$ws.Range("I126").Value = 4434.7896
$ws.Range("J126").Value = 16170.857
$ws.Range("K126").Value = 13304.3688
$ws.Range("L126").Value = 48512.571
$ws.Range("M126").Value = -10834.3688
$ws.Range("N126").Value = -53452.571
$ws.Range("H136").Value = 6132.143
$ws.Range("I136").Value = 2810
$ws.Range("K136").Value = 8430
$ws.Range("M136").Value = -5880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19931.666
$ws.Range("J41").Value = 19931.666
$ws.Range("L41").Value = 19931.666
$ws.Range("N41").Value = -20711.666
$ws.Range("H81").Value = 3823.5625
$ws.Range("I81").Value = 2811.6667
$ws.Range("K81").Value = 5623.3334
$ws.Range("M81").Value = -4562.3334
$ws.Range("H84").Value = 3823.5625
$ws.Range("I84").Value = 2811.6667
$ws.Range("K84").Value = 28116.667
$ws.Range("M84").Value = -22812.667
$ws.Range("H96").Value = 3374.5
$ws.Range("J96").Value = 3374.5
$ws.Range("L96").Value = 3374.5
$ws.Range("N96").Value = -6120.5
$ws.Range("H107").Value = 1135.7142
$ws.Range("I107").Value = 939.5
$ws.Range("K107").Value = 2818.5
$ws.Range("M107").Value = -898.5
$ws.Range("H126").Value = 2264.5217
$ws.Range("I126").Value = 2308.7273
$ws.Range("K126").Value = 6926.1819
$ws.Range("M126").Value = -4456.1819
$ws.Range("H132").Value = 3962.697
$ws.Range("I132").Value = 1665.1041
$ws.Range("J132").Value = 10089.611
$ws.Range("K132").Value = 4995.3123
$ws.Range("L132").Value = 30268.833
$ws.Range("M132").Value = -2465.3123
$ws.Range("N132").Value = -35328.833
$ws.Range("H136").Value = 2919.75
$ws.Range("I136").Value = 1682.6154
$ws.Range("K136").Value = 5047.8462
$ws.Range("M136").Value = -2497.8462
